$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp string (A1) ---
$ws.Range("A1").Value = "Datos actualizados a 13 de Septiembre de 2020 a las 12:09"

# --- Swap Bielorrusia / Polonia rows (data stays with each country) ---
$ws.Range("A48").Value = "Polonia"
$ws.Range("A49").Value = "Bielorrusia"

# --- Row 18 (Banglades) ---
$ws.Range("B18").Value = 337520
$ws.Range("C18").Value = 1476
$ws.Range("D18").Value = 240643
$ws.Range("E18").Value = 92144
$ws.Range("G18").Value = 31
$ws.Range("H18").Value = 4733

# --- Row 48 (now Polonia) ---
$ws.Range("B48").Value = 74152
$ws.Range("C48").Value = 502
$ws.Range("D48").Value = 60659
$ws.Range("E48").Value = 11305
$ws.Range("G48").Value = 6
$ws.Range("H48").Value = 2188

# --- Row 49 (now Bielorrusia) ---
$ws.Range("B49").Value = 73975
$ws.Range("D49").Value = 72547
$ws.Range("E49").Value = 684
$ws.Range("H49").Value = 744

# --- Row 76 (Australia) ---
$ws.Range("D76").Value = 23462
$ws.Range("E76").Value = 2379

# --- Row 95 ---
$ws.Range("B95").Value = 10390
$ws.Range("C95").Value = 5
$ws.Range("D95").Value = 9756
$ws.Range("E95").Value = 370
$ws.Range("G95").Value = 2
$ws.Range("H95").Value = 264

# --- Row 97 (Malasia) ---
$ws.Range("B97").Value = 9915
$ws.Range("C97").Value = 47
$ws.Range("D97").Value = 9196
$ws.Range("E97").Value = 591

# --- Row 103 (Finlandia) ---
$ws.Range("B103").Value = 8580
$ws.Range("C103").Value = 23
$ws.Range("E103").Value = 743

# --- Row 133 (Lituania) ---
$ws.Range("E133").Value = 1178
$ws.Range("G133").Value = 1
$ws.Range("H133").Value = 87

# --- Row 134 (Sri Lanka) ---
$ws.Range("B134").Value = 3204
$ws.Range("C134").Value = 9
$ws.Range("D134").Value = 2996
$ws.Range("E134").Value = 196

# --- Row 142 (Reunion) ---
$ws.Range("E142").Value = 1395
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 15
